$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed Price (D) / Volume(1h) (E) figures for the cryptos list.
# Values are written as text (NumberFormat "@") to match the source data,
# which stores prices like "66.837.93" / "1.00" / "0.0891" as literal strings
# rather than numbers, then the cell style is reset to "Normal" so no stray
# style index is left behind on the cell.
$updates = @(
    @{ Cell = "D2"; Value = '66.837.93' }
    @{ Cell = "E2"; Value = '  -3.37%  ' }
    @{ Cell = "D3"; Value = '3.469.11' }
    @{ Cell = "E3"; Value = '  -3.31%  ' }
    @{ Cell = "E4"; Value = '  +0.23%  ' }
    @{ Cell = "D5"; Value = '604.33' }
    @{ Cell = "E5"; Value = '  -3.28%  ' }
    @{ Cell = "D6"; Value = '148.31' }
    @{ Cell = "E6"; Value = '  -5.87%  ' }
    @{ Cell = "D7"; Value = '3.465.59' }
    @{ Cell = "E7"; Value = '  -3.22%  ' }
    @{ Cell = "E8"; Value = '  +0.03%  ' }
    @{ Cell = "D9"; Value = '0.483' }
    @{ Cell = "E9"; Value = '  -1.99%  ' }
    @{ Cell = "D10"; Value = '0.142' }
    @{ Cell = "E10"; Value = '  -4.34%  ' }
    @{ Cell = "D11"; Value = '7.56' }
    @{ Cell = "E11"; Value = '  +2.07%  ' }
    @{ Cell = "E12"; Value = '  -3.45%  ' }
    @{ Cell = "E13"; Value = '  -4.80%  ' }
    @{ Cell = "D14"; Value = '31.85' }
    @{ Cell = "E14"; Value = '  -5.11%  ' }
    @{ Cell = "D15"; Value = '4.059.19' }
    @{ Cell = "E15"; Value = '  -3.26%  ' }
    @{ Cell = "D16"; Value = '3.469.65' }
    @{ Cell = "E16"; Value = '  -3.31%  ' }
    @{ Cell = "D17"; Value = '66.892.28' }
    @{ Cell = "E17"; Value = '  -3.77%  ' }
    @{ Cell = "E18"; Value = '  -0.51%  ' }
    @{ Cell = "D19"; Value = '6.47' }
    @{ Cell = "E19"; Value = '  -4.48%  ' }
    @{ Cell = "D20"; Value = '15.42' }
    @{ Cell = "E20"; Value = '  -4.31%  ' }
    @{ Cell = "D21"; Value = '10.14' }
    @{ Cell = "E21"; Value = '  -0.46%  ' }
    @{ Cell = "D22"; Value = '440.69' }
    @{ Cell = "E22"; Value = '  -4.65%  ' }
    @{ Cell = "E23"; Value = '  -4.90%  ' }
    @{ Cell = "D24"; Value = '79.10' }
    @{ Cell = "E24"; Value = '  +0.34%  ' }
    @{ Cell = "E25"; Value = '  +0.14%  ' }
    @{ Cell = "D26"; Value = '3.608.38' }
    @{ Cell = "E26"; Value = '  -3.37%  ' }
    @{ Cell = "E27"; Value = '  -10.32%  ' }
    @{ Cell = "D28"; Value = '9.84' }
    @{ Cell = "E28"; Value = '  -7.89%  ' }
    @{ Cell = "D29"; Value = '8.43' }
    @{ Cell = "E29"; Value = '  -8.39%  ' }
    @{ Cell = "E30"; Value = '  -5.39%  ' }
    @{ Cell = "D31"; Value = '1.60' }
    @{ Cell = "E31"; Value = '  -6.50%  ' }
    @{ Cell = "E32"; Value = '  -2.86%  ' }
    @{ Cell = "E33"; Value = '  -0.01%  ' }
    @{ Cell = "D34"; Value = '25.44' }
    @{ Cell = "E34"; Value = '  -3.86%  ' }
    @{ Cell = "E35"; Value = '  -6.69%  ' }
    @{ Cell = "D36"; Value = '3.462.03' }
    @{ Cell = "E36"; Value = '  -3.39%  ' }
    @{ Cell = "E37"; Value = '  -6.76%  ' }
    @{ Cell = "D38"; Value = '7.95' }
    @{ Cell = "E38"; Value = '  -5.48%  ' }
    @{ Cell = "E39"; Value = '  -0.01%  ' }
    @{ Cell = "D40"; Value = '1.00' }
    @{ Cell = "E40"; Value = '  +0.19%  ' }
    @{ Cell = "D41"; Value = '174.62' }
    @{ Cell = "E41"; Value = '  -2.73%  ' }
    @{ Cell = "D42"; Value = '0.0891' }
    @{ Cell = "E42"; Value = '  -3.52%  ' }
    @{ Cell = "E43"; Value = '  -11.24%  ' }
    @{ Cell = "E44"; Value = '  -4.90%  ' }
    @{ Cell = "D45"; Value = '0.886' }
    @{ Cell = "E45"; Value = '  -2.41%  ' }
    @{ Cell = "D46"; Value = '29.48' }
    @{ Cell = "E46"; Value = '  -5.99%  ' }
    @{ Cell = "D47"; Value = '46.16' }
    @{ Cell = "E47"; Value = '  +0.32%  ' }
    @{ Cell = "D48"; Value = '1.24' }
    @{ Cell = "E48"; Value = '  -9.10%  ' }
    @{ Cell = "D49"; Value = '2.47' }
    @{ Cell = "E49"; Value = '  -9.55%  ' }
    @{ Cell = "D50"; Value = '7.48' }
    @{ Cell = "E50"; Value = '  -4.45%  ' }
    @{ Cell = "D51"; Value = '0.990' }
    @{ Cell = "E51"; Value = '  -4.78%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
